$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the slightly-off timestamp recorded for row 3 (A3)
$ws.Range("A3").Value = 45865.08354313658

# Append the new row of sensor data captured by the scheduled task (row 4)
$ws.Range("A4").Value = 45865.12527226665
$ws.Range("B4").Value = 2025
$ws.Range("C4").Value = 30
$ws.Range("D4").Value = 13.11
$ws.Range("E4").Value = 91.16
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 1.91
$ws.Range("H4").Value = "ESE"
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = "03:00:23"

# Match the date-time number format used by the other rows in column A
$ws.Range("A4").NumberFormat = $ws.Range("A3").NumberFormat
